$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Logboek update: add "Week 11" and "Week 12" entries (searchview werk).
#
# Before:
#   Row26: A26="Week 9" (wk header)           B26=<empty>
#   Row30: A30=43575 (date)  B30="1 uur 50 minuten"  C30="route beschrijvingen..."
#   Row31: (only Q filled)
#   Row32: (only Q filled)
#   Row33: (only Q filled)
#
# After:
#   Row26: B26 filled in with "10 uur 25 minuten"
#   Row30: becomes the "Week 11" header row (A30/B30)
#   Row31: gets what used to be in row30 (date 43575 entry)
#   Row32: new date entry (43579)
#   Row33: becomes the "Week 12" header row (A33 only)
#   Row34: new date entry (43586) - brand new row
# ---------------------------------------------------------------------------

# First, push the current row-30 entry (20-Apr / 43575) down into row 31,
# copying both value and number format/font so it keeps looking like the
# other date rows above it.
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B30").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A31").Value = $ws.Range("A30").Value2
$ws.Range("B31").Value = $ws.Range("B30").Value2
$ws.Range("C31").Value = $ws.Range("C30").Value2

# --- Now lay down the new content, in the same order the author typed it ---
# (this keeps newly-created shared-string entries in a sensible order)

# Row 30 becomes the "Week 11" header (style copied from the Week 9/Week 10
# header rows above).
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A30").Value = "Week 11"

# Row 26 gets its hours total filled in (style was already set on that cell).
$ws.Range("B26").Value = "10 uur 25 minuten"

# Row 32: new entry dated 43579 ("2 uur" / searchview implementeren en debuggen)
$ws.Range("A27").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("B27").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("C27").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A32").Value = 43579
$ws.Range("B32").Value = "2 uur"
$ws.Range("C32").Value = "searchview implementeren en debuggen"

# Row 33 becomes the "Week 12" header (A only).
$ws.Range("A24").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A33").Value = "Week 12"

# Row 30's hours total (same header styling as A30/B26).
$ws.Range("B24").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("B30").Value = "3 uur 50 minuten"

# The old row-30 description (column C) no longer applies to the header row -
# it has already been preserved on row 31, so clear it here.
$ws.Range("C30").ClearContents()

# Row 34: brand-new row dated 43586.
$ws.Range("A32").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A34").Value = 43586
$ws.Range("C34").Value = "searchview debuggen, onlocationChanged debuggen en oplossen, geopoints van de locatie doorsturen naar de firestore en uit de firestore uithalen in proefproject"
$ws.Range("B34").Value = "2 uur 30 minuten"

# Selection, to match what the saved file recorded.
$ws.Range("B34").Select()
